{"js": "// Goal: split the first paragraph's run \"This is a Microsoft word document.\"\n// into four runs by appending three new runs to the end of the paragraph:\n//   \" (\", \"Changed main\", \")\"\n// insertOoxml is used (instead of insertText) so the newly inserted text\n// lands in its own separate <w:r> runs rather than being silently merged\n// into the existing adjacent run.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\nconst ooxmlFragment =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>' +\n  '<w:r><w:t>Changed main</w:t></w:r>' +\n  '<w:r><w:t>)</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nconst endRange = firstParagraph.getRange(\"End\");\nendRange.insertOoxml(ooxmlFragment, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Goal: split the first paragraph's run \"This is a Microsoft word document.\"\n# into four runs by appending three new runs at the end of the paragraph:\n#   \" (\", \"Changed main\", \")\"\n#\n# Plain Range.InsertAfter() / Range.Text= calls get silently re-merged into\n# the existing adjacent run because the inserted text has identical (empty)\n# run formatting, so instead we rewrite the paragraph via Range.InsertXML()\n# with an explicit OOXML fragment that already contains the desired four\n# separate <w:r> runs. We read back the paragraph's own identity attributes\n# (paraId/textId/rsid) and its current text straight from the document first,\n# so the <w:p> element and the original run's text are reproduced unchanged,\n# exactly like the target diff (only new sibling runs are added).\n\n$d = $word.ActiveDocument\n$paragraph = $d.Paragraphs(1)\n$paragraphRange = $paragraph.Range\n\nfunction Escape-XmlText([string]$s) {\n    return $s.Replace('&', '&amp;').Replace('<', '&lt;').Replace('>', '&gt;')\n}\n\nfunction Make-RunXml([string]$text) {\n    $preserve = ''\n    if ($text.Length -gt 0 -and (($text.Substring(0,1) -eq ' ') -or ($text.Substring($text.Length-1,1) -eq ' '))) {\n        $preserve = ' xml:space=\"preserve\"'\n    }\n    return '<w:r><w:t' + $preserve + '>' + (Escape-XmlText $text) + '</w:t></w:r>'\n}\n\n# Recover the paragraph's own <w:p ...> opening tag (with its paraId/textId/\n# rsid attributes) so it is carried over untouched.\n$existingXml = $paragraphRange.WordOpenXML\n$openTag = '<w:p>'\nif ($existingXml -match '<w:p(?: [^>]*)?>') {\n    $openTag = $matches[0]\n}\n\n# Paragraph.Range.Text includes the trailing paragraph-mark character(s);\n# strip those so we only keep the visible text.\n$existingText = $paragraphRange.Text.TrimEnd([char]13, [char]7)\n\n$run0 = Make-RunXml $existingText\n$run1 = Make-RunXml \" (\"\n$run2 = Make-RunXml \"Changed main\"\n$run3 = Make-RunXml \")\"\n\n$newParagraphXml = $openTag + $run0 + $run1 + $run2 + $run3 + '</w:p>'\n\n$ooxmlPackage = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n    '<w:body>' + $newParagraphXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n$paragraphRange.InsertXML($ooxmlPackage)\n"}
